# Update column G ("K") values on Sheet1 to reflect recalculated strike
# counts (commit: "regen save_data to use K instead of Strike#, regen
# std/mean, calc and write s_vals").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> new value for column G
$gValues = @{
    2  = 2
    3  = 0
    5  = 0
    6  = 1
    7  = 0
    8  = 0
    9  = 2
    10 = 2
    11 = 1
    12 = 1
    13 = 2
    15 = 0
    16 = 1
    17 = 5
    18 = 0
    19 = 2
    20 = 1
    21 = 1
}

foreach ($row in $gValues.Keys) {
    $ws.Range("G$row").Value = $gValues[$row]
}
